$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 3932.9333
$ws.Range("I15").Value = 3932.9333
$ws.Range("K15").Value = 11798.7999
$ws.Range("M15").Value = -11629.7999

$ws.Range("H100").Value = 1392.5
$ws.Range("I100").Value = 1391
$ws.Range("J100").Value = 1400
$ws.Range("K100").Value = 1391
$ws.Range("L100").Value = 1400
$ws.Range("M100").Value = -850
$ws.Range("N100").Value = -2482

$ws.Range("H106").Value = 3691.8696
$ws.Range("I106").Value = 3405.1365
$ws.Range("K106").Value = 3405.1365
$ws.Range("M106").Value = -2774.1365

$ws.Range("H137").Value = 4168651
$ws.Range("I137").Value = 8334531.5
$ws.Range("J137").Value = 2770.25
$ws.Range("K137").Value = 25003594.5
$ws.Range("L137").Value = 8310.75
$ws.Range("M137").Value = -25001044.5
$ws.Range("N137").Value = -13410.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1841.909
$ws.Range("I2").Value = 1857.625
$ws.Range("J2").Value = 1800
$ws.Range("K2").Value = 1857.625
$ws.Range("L2").Value = 1800
$ws.Range("M2").Value = -1744.625
$ws.Range("N2").Value = -2026

$ws.Range("H32").Value = 16513.82
$ws.Range("I32").Value = 14543.436
$ws.Range("J32").Value = 30485.637
$ws.Range("K32").Value = 14543.436
$ws.Range("L32").Value = 30485.637
$ws.Range("M32").Value = -14256.436
$ws.Range("N32").Value = -31059.637

$ws.Range("H61").Value = 250250830
$ws.Range("I61").Value = 166834430
$ws.Range("K61").Value = 166834430
$ws.Range("M61").Value = -166834218

$ws.Range("H64").Value = 22091
$ws.Range("J64").Value = 22091
$ws.Range("L64").Value = 22091
$ws.Range("N64").Value = -22587

$ws.Range("H67").Value = 22091
$ws.Range("J67").Value = 22091
$ws.Range("L67").Value = 22091
$ws.Range("N67").Value = -23807

$ws.Range("H74").Value = 21001138
$ws.Range("I74").Value = 31376176
$ws.Range("J74").Value = 251058.25
$ws.Range("K74").Value = 31376176
$ws.Range("L74").Value = 251058.25
$ws.Range("M74").Value = -31375302
$ws.Range("N74").Value = -252806.25

$ws.Range("H77").Value = 21001138
$ws.Range("I77").Value = 31376176
$ws.Range("J77").Value = 251058.25
$ws.Range("K77").Value = 156880880
$ws.Range("L77").Value = 1255291.25
$ws.Range("M77").Value = -156876512
$ws.Range("N77").Value = -1264027.25

$ws.Range("H102").Value = 23811084
$ws.Range("I102").Value = 35715636
$ws.Range("J102").Value = 1979
$ws.Range("K102").Value = 35715636
$ws.Range("L102").Value = 1979
$ws.Range("M102").Value = -35714014
$ws.Range("N102").Value = -5223

$ws.Range("H116").Value = 1841.909
$ws.Range("I116").Value = 1857.625
$ws.Range("J116").Value = 1800
$ws.Range("K116").Value = 1857.625
$ws.Range("L116").Value = 1800
$ws.Range("M116").Value = 436.375
$ws.Range("N116").Value = -6388

$ws.Range("H122").Value = 2049.889
$ws.Range("I122").Value = 2039
$ws.Range("J122").Value = 2071.6667
$ws.Range("K122").Value = 6117
$ws.Range("L122").Value = 6215.000100000001
$ws.Range("M122").Value = -3667
$ws.Range("N122").Value = -11115.0001

$ws.Range("H132").Value = 102586.3
$ws.Range("I132").Value = 92611.45
$ws.Range("J132").Value = 114777.78
$ws.Range("K132").Value = 277834.35
$ws.Range("L132").Value = 344333.34
$ws.Range("M132").Value = -275304.35
$ws.Range("N132").Value = -349393.34

$ws.Range("H136").Value = 250250830
$ws.Range("I136").Value = 166834430
$ws.Range("K136").Value = 500503290
$ws.Range("M136").Value = -500500740

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1841.909
$ws.Range("I3").Value = 1857.625
$ws.Range("J3").Value = 1800
$ws.Range("K3").Value = 1857.625
$ws.Range("L3").Value = 1800
$ws.Range("M3").Value = -1743.625
$ws.Range("N3").Value = -2028

$ws.Range("H134").Value = 3773.8333
$ws.Range("I134").Value = 4668.25
$ws.Range("J134").Value = 1985
$ws.Range("K134").Value = 14004.75
$ws.Range("L134").Value = 5955
$ws.Range("M134").Value = -11469.75
$ws.Range("N134").Value = -11025

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2540.7778
$ws.Range("I31").Value = 1005.20514
$ws.Range("J31").Value = 6533.2666
$ws.Range("K31").Value = 1005.20514
$ws.Range("L31").Value = 6533.2666
$ws.Range("M31").Value = -710.20514
$ws.Range("N31").Value = -7123.2666

$ws.Range("H34").Value = 2540.7778
$ws.Range("I34").Value = 1005.20514
$ws.Range("J34").Value = 6533.2666
$ws.Range("K34").Value = 1005.20514
$ws.Range("L34").Value = 6533.2666
$ws.Range("M34").Value = -803.20514
$ws.Range("N34").Value = -6937.2666

$ws.Range("H58").Value = 46513372
$ws.Range("I58").Value = 45455840
$ws.Range("J58").Value = 47621260
$ws.Range("K58").Value = 45455840
$ws.Range("L58").Value = 47621260
$ws.Range("M58").Value = -45455637
$ws.Range("N58").Value = -47621666

$ws.Range("H132").Value = 45180.78
$ws.Range("I132").Value = 1549.1177
$ws.Range("K132").Value = 4647.3531
$ws.Range("M132").Value = -2117.3531

$ws.Range("H134").Value = 36639.418
$ws.Range("I134").Value = 1825.24
$ws.Range("J134").Value = 181698.5
$ws.Range("K134").Value = 5475.72
$ws.Range("L134").Value = 545095.5
$ws.Range("M134").Value = -2940.72
$ws.Range("N134").Value = -550165.5

$ws.Range("H136").Value = 46513372
$ws.Range("I136").Value = 45455840
$ws.Range("J136").Value = 47621260
$ws.Range("K136").Value = 136367520
$ws.Range("L136").Value = 142863780
$ws.Range("M136").Value = -136364970
$ws.Range("N136").Value = -142868880

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1450.8334
$ws.Range("I132").Value = 858.5714
$ws.Range("J132").Value = 1827.7273
$ws.Range("K132").Value = 7727.1426
$ws.Range("L132").Value = 16449.5457
$ws.Range("M132").Value = -5197.1426
$ws.Range("N132").Value = -21509.5457

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 225544.33
$ws.Range("I132").Value = 1000000
$ws.Range("J132").Value = 128737.375
$ws.Range("K132").Value = 3000000
$ws.Range("L132").Value = 386212.125
$ws.Range("M132").Value = -2997470
$ws.Range("N132").Value = -391272.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H34").Value = 70024
$ws.Range("J34").Value = 70024
$ws.Range("L34").Value = 70024
$ws.Range("N34").Value = -70368

$ws.Range("H40").Value = 4000
$ws.Range("I40").Value = 3400
$ws.Range("J40").Value = 5500
$ws.Range("K40").Value = 3400
$ws.Range("L40").Value = 5500
$ws.Range("M40").Value = -3264
$ws.Range("N40").Value = -5772

$ws.Range("H55").Value = 105.4
$ws.Range("I55").Value = 62.5
$ws.Range("J55").Value = 277
$ws.Range("K55").Value = 62.5
$ws.Range("L55").Value = 277
$ws.Range("M55").Value = 110.5
$ws.Range("N55").Value = -623

$ws.Range("H61").Value = 1717.6818
$ws.Range("I61").Value = 1766.4706
$ws.Range("J61").Value = 1551.8
$ws.Range("K61").Value = 1766.4706
$ws.Range("L61").Value = 1551.8
$ws.Range("M61").Value = -1564.4706
$ws.Range("N61").Value = -1955.8

$ws.Range("H93").Value = 2851.3333
$ws.Range("J93").Value = 2851.3333
$ws.Range("L93").Value = 2851.3333
$ws.Range("N93").Value = -5347.3333

$ws.Range("H113").Value = 1717.6818
$ws.Range("I113").Value = 1766.4706
$ws.Range("J113").Value = 1551.8
$ws.Range("K113").Value = 1766.4706
$ws.Range("L113").Value = 1551.8
$ws.Range("M113").Value = 403.5293999999999
$ws.Range("N113").Value = -5891.8

$ws.Range("H132").Value = 94990.73
$ws.Range("I132").Value = 2475
$ws.Range("J132").Value = 147856.86
$ws.Range("K132").Value = 7425
$ws.Range("L132").Value = 443570.58
$ws.Range("M132").Value = -4895
$ws.Range("N132").Value = -448630.58

$ws.Range("H136").Value = 183143.9
$ws.Range("I136").Value = 251076
$ws.Range("J136").Value = 144325.58
$ws.Range("K136").Value = 753228
$ws.Range("L136").Value = 432976.74
$ws.Range("M136").Value = -750678
$ws.Range("N136").Value = -438076.74

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 101770.2
$ws.Range("I132").Value = 84458
$ws.Range("K132").Value = 253374
$ws.Range("M132").Value = -250844

$ws.Range("H136").Value = 44149.914
$ws.Range("I136").Value = 26568.871
$ws.Range("J136").Value = 129857.5
$ws.Range("K136").Value = 79706.613
$ws.Range("L136").Value = 389572.5
$ws.Range("M136").Value = -77156.613
$ws.Range("N136").Value = -394672.5
